# Update "Pagos" (F) and "Inscrições homologadas" (H) columns for the
# rows whose registration counts increased by one payment each.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(29, 33, 36, 38, 55, 79)

foreach ($r in $rows) {
    $ws.Range("F$r").Value2 = $ws.Range("F$r").Value2 + 1
    $ws.Range("H$r").Value2 = $ws.Range("H$r").Value2 + 1
}
